# Calibration edits: update the New-This-Year share values on the
# passenger (psgr) and freight (frgt) sheets, and restore the selection
# state left behind on each sheet after the calibration pass.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# SoCDTtiNTY-psgr (passenger)
# ---------------------------------------------------------------------
$psgr = $wb.Worksheets.Item("SoCDTtiNTY-psgr")

$psgr.Range("B2:H2").Value = 0.060966558335291694   # LDVs
$psgr.Range("B3:H3").Value = 0.1                    # HDVs
$psgr.Range("B4:H4").Value = 0.05                   # aircraft
$psgr.Range("B5:H5").Value = 0.028571                # rail
$psgr.Range("B6:H6").Value = 0.030303                # ships
$psgr.Range("B7:H7").Value = 0.1                    # motorbikes

$psgr.Activate()
$psgr.Range("B2:C2").Select() | Out-Null


# ---------------------------------------------------------------------
# SoCDTtiNTY-frgt (freight)
# ---------------------------------------------------------------------
$frgt = $wb.Worksheets.Item("SoCDTtiNTY-frgt")

$frgt.Range("B2:H2").Value = 0.1                    # LDVs
$frgt.Range("B3:H3").Value = 0.1                    # HDVs
$frgt.Range("B4:H4").Value = 0.03125                 # aircraft
$frgt.Range("B5:H5").Value = 0.028571                # rail
$frgt.Range("B6:H6").Value = 0.030303                # ships
$frgt.Range("B7:H7").Value = 0.1                    # motorbikes

$frgt.Activate()
$frgt.Range("B2:H7").Select() | Out-Null

# The "About" sheet was (and remains) the selected tab in the source file.
$about = $wb.Worksheets.Item("About")
$about.Activate()
